$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) on this sheet stores plain-text values (some look like
# numbers with multiple "." separators, e.g. "1.639.54", which are not valid
# numeric values). Force NumberFormat to Text on each touched D cell before
# writing so Excel does not auto-convert the string into a floating point number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.709.92"
$ws.Range("E2").Value = "  +0.57%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.639.95"
$ws.Range("E3").Value = "  -0.48%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.67"
$ws.Range("E5").Value = "  +0.03%  "

$ws.Range("E6").Value = "  -2.00%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.23"
$ws.Range("E8").Value = "  -1.39%  "

$ws.Range("E9").Value = "  +1.08%  "

$ws.Range("E10").Value = "  +0.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0889"
$ws.Range("E11").Value = "  -0.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.871.60"
$ws.Range("E12").Value = "  -0.49%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.628.57"
$ws.Range("E13").Value = "  -1.15%  "

$ws.Range("E14").Value = "  +0.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.563"
$ws.Range("E15").Value = "  -4.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.79"
$ws.Range("E16").Value = "  +0.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.669.95"
$ws.Range("E17").Value = "  +0.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.73"
$ws.Range("E18").Value = "  +0.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.71"
$ws.Range("E19").Value = "  +2.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0723"
$ws.Range("E20").Value = "  -0.30%  "

$ws.Range("E21").Value = "  +0.00%  "

$ws.Range("E22").Value = "  -0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.28"
$ws.Range("E23").Value = "  +5.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.06"
$ws.Range("E24").Value = "  +2.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.86"
$ws.Range("E25").Value = "  +1.48%  "

$ws.Range("E26").Value = "  -1.01%  "

$ws.Range("E27").Value = "  -0.71%  "

$ws.Range("E28").Value = "  -0.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.61"
$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("E30").Value = "  +0.19%  "

$ws.Range("E31").Value = "  -0.10%  "

$ws.Range("E32").Value = "  +0.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.458.83"
$ws.Range("E33").Value = "  +2.27%  "

$ws.Range("E34").Value = "  -2.26%  "

$ws.Range("E35").Value = "  -2.16%  "

$ws.Range("E36").Value = "  -0.34%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.569"
$ws.Range("E37").Value = "  +0.12%  "

$ws.Range("E38").Value = "  -0.43%  "

$ws.Range("E39").Value = "  +0.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.897"
$ws.Range("E40").Value = "  +9.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.31"
$ws.Range("E41").Value = "  +6.33%  "

$ws.Range("E43").Value = "  -1.68%  "

$ws.Range("E44").Value = "  +0.82%  "

$ws.Range("E45").Value = "  -0.82%  "

$ws.Range("E46").Value = "  -0.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.781.47"
$ws.Range("E47").Value = "  -0.54%  "

$ws.Range("E48").Value = "  +3.25%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.02"
$ws.Range("E49").Value = "  -1.25%  "

$ws.Range("E50").Value = "  -1.34%  "

$ws.Range("E51").Value = "  +0.02%  "

